$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:R2").Formula = '=2*ABS(C$1-$A2)^2'
$ws.Range("C3:R17").Formula = '=ABS(C$1-$A3)^2'
$ws.Range("N26").Select() | Out-Null
